$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleD2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.931.90"
$ws.Range("D2").Style = $styleD2
$ws.Range("E2").Value = "  -1.52%  "
$styleD3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.14"
$ws.Range("D3").Style = $styleD3
$ws.Range("E3").Value = "  +0.06%  "
$styleD4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9955"
$ws.Range("D4").Style = $styleD4
$ws.Range("E4").Value = "  +0.51%  "
$styleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.68"
$ws.Range("D5").Style = $styleD5
$ws.Range("E5").Value = "  -2.96%  "
$styleD6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9956"
$ws.Range("D6").Style = $styleD6
$ws.Range("E6").Value = "  +0.33%  "
$styleD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4937"
$ws.Range("D7").Style = $styleD7
$ws.Range("E7").Value = "  -0.81%  "
$styleD8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.17"
$ws.Range("D8").Style = $styleD8
$ws.Range("E8").Value = "  -1.22%  "
$styleD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2901"
$ws.Range("D9").Style = $styleD9
$ws.Range("E9").Value = "  +1.60%  "
$styleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06606"
$ws.Range("D10").Style = $styleD10
$ws.Range("E10").Value = "  +0.66%  "
$styleD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.874.17"
$ws.Range("D11").Style = $styleD11
$ws.Range("E11").Value = "  +1.44%  "
$styleD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.91"
$ws.Range("D12").Style = $styleD12
$ws.Range("E12").Value = "  -1.14%  "
$styleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07171"
$ws.Range("D13").Style = $styleD13
$ws.Range("E13").Value = "  +0.02%  "
$styleD14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6692"
$ws.Range("D14").Style = $styleD14
$ws.Range("E14").Value = "  -0.30%  "
$styleD15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.62"
$ws.Range("D15").Style = $styleD15
$ws.Range("E15").Value = "  -1.06%  "
$styleD16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.809"
$ws.Range("D16").Style = $styleD16
$ws.Range("E16").Value = "  +0.22%  "
$styleD17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.904.79"
$ws.Range("D17").Style = $styleD17
$ws.Range("E17").Value = "  -1.22%  "
$styleD18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007875"
$ws.Range("D18").Style = $styleD18
$ws.Range("E18").Value = "  +4.86%  "
$styleD19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9975"
$ws.Range("D19").Style = $styleD19
$ws.Range("E19").Value = "  +0.34%  "
$styleD20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("D20").Style = $styleD20
$ws.Range("E20").Value = "  +1.20%  "
$styleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.66"
$ws.Range("D21").Style = $styleD21
$ws.Range("E21").Value = "  +1.63%  "
$styleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9949"
$ws.Range("D22").Style = $styleD22
$ws.Range("E22").Value = "  +0.65%  "
$styleD23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.748"
$ws.Range("D23").Style = $styleD23
$ws.Range("E23").Value = "  +1.36%  "
$styleD24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.576"
$ws.Range("D24").Style = $styleD24
$ws.Range("E24").Value = "  +1.17%  "
$styleD25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.087"
$ws.Range("D25").Style = $styleD25
$ws.Range("E25").Value = "  +1.04%  "
$styleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.99"
$ws.Range("D26").Style = $styleD26
$ws.Range("E26").Value = "  +2.10%  "
$styleD27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.39"
$ws.Range("D27").Style = $styleD27
$ws.Range("E27").Value = "  -1.10%  "
$styleD28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.70"
$ws.Range("D28").Style = $styleD28
$ws.Range("E28").Value = "  -0.72%  "
$styleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.925"
$ws.Range("D29").Style = $styleD29
$ws.Range("E29").Value = "  -1.59%  "
$styleD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.377"
$ws.Range("D30").Style = $styleD30
$ws.Range("E30").Value = "  -2.51%  "
$styleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.177"
$ws.Range("D31").Style = $styleD31
$ws.Range("E31").Value = "  -1.77%  "
$styleD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08598"
$ws.Range("D32").Style = $styleD32
$ws.Range("E32").Value = "  -0.30%  "
$styleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.922"
$ws.Range("D33").Style = $styleD33
$ws.Range("E33").Value = "  +0.60%  "
$styleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04983"
$ws.Range("D34").Style = $styleD34
$ws.Range("E34").Value = "  -1.26%  "
$styleD35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("D35").Style = $styleD35
$ws.Range("E35").Value = "  -2.26%  "
$styleD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7048"
$ws.Range("D36").Style = $styleD36
$ws.Range("E36").Value = "  +2.60%  "
$styleD37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.650"
$ws.Range("D37").Style = $styleD37
$ws.Range("E37").Value = "  -1.33%  "
$styleD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.204"
$ws.Range("D38").Style = $styleD38
$ws.Range("E38").Value = "  -5.37%  "
$styleD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.679"
$ws.Range("D39").Style = $styleD39
$ws.Range("E39").Value = "  -3.02%  "
$styleD40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9318"
$ws.Range("D40").Style = $styleD40
$ws.Range("E40").Value = "  -2.17%  "
$styleD41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01635"
$ws.Range("D41").Style = $styleD41
$ws.Range("E41").Value = "  +0.41%  "
$styleD42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.059"
$ws.Range("D42").Style = $styleD42
$ws.Range("E42").Value = "  -1.17%  "
$styleD43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9931"
$ws.Range("D43").Style = $styleD43
$ws.Range("E43").Value = "  -0.06%  "
$styleD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.47"
$ws.Range("D44").Style = $styleD44
$ws.Range("E44").Value = "  -1.15%  "
$styleD45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4178"
$ws.Range("D45").Style = $styleD45
$ws.Range("E45").Value = "  +0.25%  "
$styleD46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.548"
$ws.Range("D46").Style = $styleD46
$ws.Range("E46").Value = "  +1.47%  "
$styleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1256"
$ws.Range("D47").Style = $styleD47
$ws.Range("E47").Value = "  +0.86%  "
$styleD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05696"
$ws.Range("D48").Style = $styleD48
$ws.Range("E48").Value = "  +1.36%  "
$styleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.53"
$ws.Range("D49").Style = $styleD49
$ws.Range("E49").Value = "  +0.39%  "
$styleD50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.227"
$ws.Range("D50").Style = $styleD50
$ws.Range("E50").Value = "  -0.83%  "
$styleD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3703"
$ws.Range("D51").Style = $styleD51
$ws.Range("E51").Value = "  -0.10%  "

Write-Host "Updated crypto prices"